# The three species-observation records currently stored in rows 19, 20
# and 21 get rotated one position:
#   new row 19 <- old row 20  (Kortskaftad ärgspik / Microcalicium ahlneri)
#   new row 20 <- old row 21  (Brunpudrad nållav / Chaenotheca gracillima)
#   new row 21 <- old row 19  (Tretåig hackspett / Picoides tridactylus)
#
# Only the columns whose content actually differs between the three
# records are written (Id, Taxonsorteringsordning, TaxonId, Artnamn,
# Vetenskapligt namn, Auktor, Enhet, Kön, Aktivitet, Ost, Nord, Publik
# kommentar, Bestämningsmetod). Every other column (C, D, I, K, N, P, S,
# T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY) holds the same
# value in all three rows, so it is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-OrClear($addr, $value) {
    if ($null -eq $value -or $value -eq "") {
        $ws.Range($addr).ClearContents()
    } else {
        $ws.Range($addr).Value2 = $value
    }
}

# --- new row 19 (was row 20's record) ---------------------------------
$ws.Range("A19").Value2 = 111742256
$ws.Range("B19").Value2 = 79444
$ws.Range("E19").Value2 = 1049
$ws.Range("F19").Value2 = "Kortskaftad ärgspik"
$ws.Range("G19").Value2 = "Microcalicium ahlneri"
$ws.Range("H19").Value2 = "Tibell"
Set-OrClear "J19" ""
Set-OrClear "L19" $null
Set-OrClear "M19" $null
$ws.Range("Q19").Value2 = 331773.1827125447
$ws.Range("R19").Value2 = 6626566.53343309
Set-OrClear "AC19" "På barklös talltorraka"
Set-OrClear "AF19" ""

# --- new row 20 (was row 21's record) ---------------------------------
$ws.Range("A20").Value2 = 111741735
$ws.Range("B20").Value2 = 73689
$ws.Range("E20").Value2 = 308
$ws.Range("F20").Value2 = "Brunpudrad nållav"
$ws.Range("G20").Value2 = "Chaenotheca gracillima"
$ws.Range("H20").Value2 = "(Vain.) Tibell"
Set-OrClear "J20" ""
Set-OrClear "L20" $null
Set-OrClear "M20" $null
$ws.Range("Q20").Value2 = 331238.0752669616
$ws.Range("R20").Value2 = 6626585.695077832
Set-OrClear "AC20" "På gran"
Set-OrClear "AF20" ""

# --- new row 21 (was row 19's record) ---------------------------------
$ws.Range("A21").Value2 = 111741430
$ws.Range("B21").Value2 = 56398
$ws.Range("E21").Value2 = 100109
$ws.Range("F21").Value2 = "Tretåig hackspett"
$ws.Range("G21").Value2 = "Picoides tridactylus"
$ws.Range("H21").Value2 = "(Linnaeus, 1758)"
Set-OrClear "J21" $null
Set-OrClear "L21" ""
Set-OrClear "M21" "färska spår"
$ws.Range("Q21").Value2 = 331285.2567537006
$ws.Range("R21").Value2 = 6626678.453820148
Set-OrClear "AC21" $null
Set-OrClear "AF21" $null
